$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), copying the header formatting used by the
# existing header cells (e.g. G1 "sum") so it picks up the same style index.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the Save column on row 2.
$ws.Range("H2").Value = 1
